$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet positioned right after the existing one
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "FTNC_Average_Demand151"
$ws2.Outline.SummaryRow = 1
$ws2.Outline.SummaryColumn = 1

# Copy the values + formatting from the original sheet (skip A1, which is empty there)
$ws1.Range("B1:F1").Copy($ws2.Range("B1"))
$ws1.Range("A2:F2").Copy($ws2.Range("A2"))

# Overwrite the figures in row 2 on the new sheet with their updated values
$ws2.Range("A2").Value = "FTNC_Average_Demand_15"
$ws2.Range("B2").Value = 2425.00693602971
$ws2.Range("C2").Value = 12847.37736604325
$ws2.Range("D2").Value = 859.3641716139318
$ws2.Range("E2").Value = 12.72209238597305
$ws2.Range("F2").Value = 16144.47056636828
